$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BurnUp Chart")

$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 28
$ws.Range("F2").Value = 31
$ws.Range("G2").Value = 39

$ws.Range("M15").Select()
